# Adds two new columns (lat_dd / long_dd) to Sheet1, matching the
# formatting of the existing last column (J) row-by-row, and fills in
# known lat/long values for the three Washington "current" rows (6-8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (K1, L1) - new shared strings "lat_dd" / "long_dd".
$ws.Range("K1").Value = "lat_dd"
$ws.Range("L1").Value = "long_dd"

# Stamp the formatting of column J (rows 2-29) onto the new K and L
# columns so each row's fill/style carries across, without touching any
# existing data in columns A-J.
$ws.Range("J2:J29").Copy()
$ws.Range("K2:K29").PasteSpecial(-4122)
$ws.Range("L2:L29").PasteSpecial(-4122)

# Known coordinates for the three newly-geocoded rows.
$ws.Range("K6").Value = 46.955267
$ws.Range("L6").Value = -124.050733

$ws.Range("K7").Value = 46.540933
$ws.Range("L7").Value = -123.972545

$ws.Range("K8").Value = 46.283661
$ws.Range("L8").Value = -123.704545

# Leave the selection where the author last left it.
$ws.Range("L5").Select()
